$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor expression values and derived edge weights/specificities
# (new TPM-based values per commit "update scripts wuth new tpm")

# Row 2
$ws.Range("G2").Value = 0.7428659999999999
$ws.Range("H2").Value = 2.228598
$ws.Range("I2").Value = 0.07647842579626549
$ws.Range("J2").Value = 0.07647842579626549
$ws.Range("M2").Value = 0.9613523333333335
$ws.Range("N2").Value = 2.884057
$ws.Range("O2").Value = 0.117240948966235
$ws.Range("P2").Value = 0.117240948966235
$ws.Range("Q2").Value = 0.714155962454
$ws.Range("R2").Value = 6.427403662086
$ws.Range("S2").Value = 0.00896640321579795
$ws.Range("T2").Value = 0.00896640321579795

# Row 3
$ws.Range("G3").Value = 0.7428659999999999
$ws.Range("H3").Value = 2.228598
$ws.Range("I3").Value = 0.07647842579626549
$ws.Range("J3").Value = 0.07647842579626549
$ws.Range("O3").Value = 0.17360760831565
$ws.Range("P3").Value = 0.17360760831565
$ws.Range("Q3").Value = 1.057505160946
$ws.Range("R3").Value = 9.517546448513999
$ws.Range("S3").Value = 0.01327723659023556
$ws.Range("T3").Value = 0.01327723659023556

# Row 4
$ws.Range("G4").Value = 0.7428659999999999
$ws.Range("H4").Value = 2.228598
$ws.Range("I4").Value = 0.07647842579626549
$ws.Range("J4").Value = 0.07647842579626549
$ws.Range("M4").Value = 1.214815
$ws.Range("N4").Value = 3.644445
$ws.Range("O4").Value = 0.1481517841898583
$ws.Range("P4").Value = 0.1481517841898583
$ws.Range("Q4").Value = 0.9024447597899998
$ws.Range("R4").Value = 8.122002838109999
$ws.Range("S4").Value = 0.01133041523374842
$ws.Range("T4").Value = 0.01133041523374842

# Row 5
$ws.Range("G5").Value = 0.7428659999999999
$ws.Range("H5").Value = 2.228598
$ws.Range("I5").Value = 0.07647842579626549
$ws.Range("J5").Value = 0.07647842579626549
$ws.Range("M5").Value = 0.3083506666666667
$ws.Range("N5").Value = 0.925052
$ws.Range("O5").Value = 0.03760465702415506
$ws.Range("P5").Value = 0.03760465702415506
$ws.Range("Q5").Value = 0.229063226344
$ws.Range("R5").Value = 2.061569037096
$ws.Range("S5").Value = 0.002875944971815857
$ws.Range("T5").Value = 0.002875944971815857

# Row 6
$ws.Range("G6").Value = 0.7428659999999999
$ws.Range("H6").Value = 2.228598
$ws.Range("I6").Value = 0.07647842579626549
$ws.Range("J6").Value = 0.07647842579626549
$ws.Range("M6").Value = 3.81316
$ws.Range("N6").Value = 11.43948
$ws.Range("O6").Value = 0.4650308544110832
$ws.Range("P6").Value = 0.4650308544110832
$ws.Range("Q6").Value = 2.83266691656
$ws.Range("R6").Value = 25.49400224904
$ws.Range("S6").Value = 0.03556482769205196
$ws.Range("T6").Value = 0.03556482769205196

# Row 7
$ws.Range("G7").Value = 0.7428659999999999
$ws.Range("H7").Value = 2.228598
$ws.Range("I7").Value = 0.07647842579626549
$ws.Range("J7").Value = 0.07647842579626549
$ws.Range("M7").Value = 0.4785743333333333
$ws.Range("N7").Value = 1.435723
$ws.Range("O7").Value = 0.05836414709301852
$ws.Range("P7").Value = 0.05836414709301852
$ws.Range("Q7").Value = 0.3555166007059999
$ws.Range("R7").Value = 3.199649406353999
$ws.Range("S7").Value = 0.004463598092615741
$ws.Range("T7").Value = 0.004463598092615741

# Row 8
$ws.Range("G8").Value = 1.185428333333333
$ws.Range("H8").Value = 3.556285
$ws.Range("I8").Value = 0.1220404390934893
$ws.Range("J8").Value = 0.1220404390934893
$ws.Range("M8").Value = 0.9613523333333335
$ws.Range("N8").Value = 2.884057
$ws.Range("O8").Value = 0.117240948966235
$ws.Range("P8").Value = 0.117240948966235
$ws.Range("Q8").Value = 1.139614294249444
$ws.Range("R8").Value = 10.256528648245
$ws.Range("S8").Value = 0.01430813689157668
$ws.Range("T8").Value = 0.01430813689157668

# Row 9
$ws.Range("G9").Value = 1.185428333333333
$ws.Range("H9").Value = 3.556285
$ws.Range("I9").Value = 0.1220404390934893
$ws.Range("J9").Value = 0.1220404390934893
$ws.Range("O9").Value = 0.17360760831565
$ws.Range("P9").Value = 0.17360760831565
$ws.Range("Q9").Value = 1.687513737917222
$ws.Range("R9").Value = 15.187623641255
$ws.Range("S9").Value = 0.02118714874881242
$ws.Range("T9").Value = 0.02118714874881242

# Row 10
$ws.Range("G10").Value = 1.185428333333333
$ws.Range("H10").Value = 3.556285
$ws.Range("I10").Value = 0.1220404390934893
$ws.Range("J10").Value = 0.1220404390934893
$ws.Range("M10").Value = 1.214815
$ws.Range("N10").Value = 3.644445
$ws.Range("O10").Value = 0.1481517841898583
$ws.Range("P10").Value = 0.1481517841898583
$ws.Range("Q10").Value = 1.440076120758333
$ws.Range("R10").Value = 12.960685086825
$ws.Range("S10").Value = 0.01808050879501417
$ws.Range("T10").Value = 0.01808050879501417

# Row 11
$ws.Range("G11").Value = 1.185428333333333
$ws.Range("H11").Value = 3.556285
$ws.Range("I11").Value = 0.1220404390934893
$ws.Range("J11").Value = 0.1220404390934893
$ws.Range("M11").Value = 0.3083506666666667
$ws.Range("N11").Value = 0.925052
$ws.Range("O11").Value = 0.03760465702415506
$ws.Range("P11").Value = 0.03760465702415506
$ws.Range("Q11").Value = 0.3655276168688888
$ws.Range("R11").Value = 3.28974855182
$ws.Range("S11").Value = 0.004589288855187949
$ws.Range("T11").Value = 0.004589288855187949

# Row 12
$ws.Range("G12").Value = 1.185428333333333
$ws.Range("H12").Value = 3.556285
$ws.Range("I12").Value = 0.1220404390934893
$ws.Range("J12").Value = 0.1220404390934893
$ws.Range("M12").Value = 3.81316
$ws.Range("N12").Value = 11.43948
$ws.Range("O12").Value = 0.4650308544110832
$ws.Range("P12").Value = 0.4650308544110832
$ws.Range("Q12").Value = 4.520227903533333
$ws.Range("R12").Value = 40.68205113179999
$ws.Range("S12").Value = 0.05675256966434909
$ws.Range("T12").Value = 0.05675256966434908

# Row 13
$ws.Range("G13").Value = 1.185428333333333
$ws.Range("H13").Value = 3.556285
$ws.Range("I13").Value = 0.1220404390934893
$ws.Range("J13").Value = 0.1220404390934893
$ws.Range("M13").Value = 0.4785743333333333
$ws.Range("N13").Value = 1.435723
$ws.Range("O13").Value = 0.05836414709301852
$ws.Range("P13").Value = 0.05836414709301852
$ws.Range("Q13").Value = 0.5673155743394444
$ws.Range("R13").Value = 5.105840169054999
$ws.Range("S13").Value = 0.007122786138548977
$ws.Range("T13").Value = 0.007122786138548976

# Row 14
$ws.Range("G14").Value = 7.785111666666666
$ws.Range("H14").Value = 23.355335
$ws.Range("I14").Value = 0.8014811351102453
$ws.Range("J14").Value = 0.8014811351102452
$ws.Range("M14").Value = 0.9613523333333335
$ws.Range("N14").Value = 2.884057
$ws.Range("O14").Value = 0.117240948966235
$ws.Range("P14").Value = 0.117240948966235
$ws.Range("Q14").Value = 7.484235266010557
$ws.Range("R14").Value = 67.358117394095
$ws.Range("S14").Value = 0.09396640885886035
$ws.Range("T14").Value = 0.09396640885886033

# Row 15
$ws.Range("G15").Value = 7.785111666666666
$ws.Range("H15").Value = 23.355335
$ws.Range("I15").Value = 0.8014811351102453
$ws.Range("J15").Value = 0.8014811351102452
$ws.Range("O15").Value = 0.17360760831565
$ws.Range("P15").Value = 0.17360760831565
$ws.Range("Q15").Value = 11.08247754782278
$ws.Range("R15").Value = 99.742297930405
$ws.Range("S15").Value = 0.139143222976602
$ws.Range("T15").Value = 0.139143222976602

# Row 16
$ws.Range("G16").Value = 7.785111666666666
$ws.Range("H16").Value = 23.355335
$ws.Range("I16").Value = 0.8014811351102453
$ws.Range("J16").Value = 0.8014811351102452
$ws.Range("M16").Value = 1.214815
$ws.Range("N16").Value = 3.644445
$ws.Range("O16").Value = 0.1481517841898583
$ws.Range("P16").Value = 0.1481517841898583
$ws.Range("Q16").Value = 9.457470429341667
$ws.Range("R16").Value = 85.117233864075
$ws.Range("S16").Value = 0.1187408601610957
$ws.Range("T16").Value = 0.1187408601610957

# Row 17
$ws.Range("G17").Value = 7.785111666666666
$ws.Range("H17").Value = 23.355335
$ws.Range("I17").Value = 0.8014811351102453
$ws.Range("J17").Value = 0.8014811351102452
$ws.Range("M17").Value = 0.3083506666666667
$ws.Range("N17").Value = 0.925052
$ws.Range("O17").Value = 0.03760465702415506
$ws.Range("P17").Value = 0.03760465702415506
$ws.Range("Q17").Value = 2.400544372491111
$ws.Range("R17").Value = 21.60489935242
$ws.Range("S17").Value = 0.03013942319715126
$ws.Range("T17").Value = 0.03013942319715125

# Row 18
$ws.Range("G18").Value = 7.785111666666666
$ws.Range("H18").Value = 23.355335
$ws.Range("I18").Value = 0.8014811351102453
$ws.Range("J18").Value = 0.8014811351102452
$ws.Range("M18").Value = 3.81316
$ws.Range("N18").Value = 11.43948
$ws.Range("O18").Value = 0.4650308544110832
$ws.Range("P18").Value = 0.4650308544110832
$ws.Range("Q18").Value = 29.68587640286666
$ws.Range("R18").Value = 267.1728876258
$ws.Range("S18").Value = 0.3727134570546822
$ws.Range("T18").Value = 0.3727134570546822

# Row 19
$ws.Range("G19").Value = 7.785111666666666
$ws.Range("H19").Value = 23.355335
$ws.Range("I19").Value = 0.8014811351102453
$ws.Range("J19").Value = 0.8014811351102452
$ws.Range("M19").Value = 0.4785743333333333
$ws.Range("N19").Value = 1.435723
$ws.Range("O19").Value = 0.05836414709301852
$ws.Range("P19").Value = 0.05836414709301852
$ws.Range("Q19").Value = 3.725754625800555
$ws.Range("R19").Value = 33.531791632205
$ws.Range("S19").Value = 0.04677776286185381
$ws.Range("T19").Value = 0.0467777628618538
